# Update "Forecast Comparison" sheet with corrected forecast output:
#  - insert a new "Week_Start_Date" column (B), shifting ASIN..is_holiday_week right by one
#  - normalize the Week labels from "W01".."W16" to "W1".."W16" (no leading zero)
#  - fill in the new Week_Start_Date column with the week's start date (as text)
#  - store is_holiday_week as a proper boolean value

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank column before the current column B (ASIN), shifting
# ASIN / MyForecast / Amazon Mean / P70 / P80 / P90 / Product Title / is_holiday_week
# one column to the right (B -> C, C -> D, ... I -> J).
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# Make sure the new date-like strings are stored as literal text, not parsed
# into Excel date serials.
$ws.Range("B2:B17").NumberFormat = "@"

$weeks = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
$weekStartDates = @("2025-01-05","2025-01-12","2025-01-19","2025-01-26","2025-02-02","2025-02-09","2025-02-16","2025-02-23","2025-03-02","2025-03-09","2025-03-16","2025-03-23","2025-03-30","2025-04-06","2025-04-13","2025-04-20")

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2

    # Week label without the leading zero (W01 -> W1, ... W10 stays W10).
    $ws.Cells.Item($row, 1).Value = $weeks[$i]

    # New Week_Start_Date column (text).
    $ws.Cells.Item($row, 2).Value = $weekStartDates[$i]

    # is_holiday_week (now column J) becomes a real boolean.
    $ws.Cells.Item($row, 10).Value = $false
}
